$wb = $excel.ActiveWorkbook

# --- Update values on the "data" sheet ---
$wsData = $wb.Worksheets.Item("data")
$wsData.Range("C2").Value = 0.04
$wsData.Range("B3").Value = 124
$wsData.Range("C3").Value = 0.04
$wsData.Range("B5").Value = 124

# --- Update selection on "columnSpecs" sheet (no longer the active tab) ---
$wsColumnSpecs = $wb.Worksheets.Item("columnSpecs")
$wsColumnSpecs.Activate()
$wsColumnSpecs.Range("G7").Select()

# --- Update selection on "data" sheet and make it the active tab ---
$wsData.Activate()
$wsData.Range("F8").Select()
